# Regenerate orders with updated distance/size codes.
# Applies global substring replacements to the textual values in the sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S20 and S25 remain unchanged)
# These codes appear (possibly combined) inside the Condition, Filename_Left,
# Filename_Right, Distance and Size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value2
    if ($val -is [string]) {
      $newVal = $val.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
      if ($newVal -ne $val) {
        $cell.Value2 = $newVal
      }
    }
  }
}
